$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 307.3793
$ws.Range("I6").Value = 131.07692
$ws.Range("J6").Value = 450.625
$ws.Range("K6").Value = 393.23076
$ws.Range("L6").Value = 1351.875
$ws.Range("M6").Value = -281.23076
$ws.Range("N6").Value = -1575.875
$ws.Range("H38").Value = 2316.8667
$ws.Range("I38").Value = 1361.5555
$ws.Range("J38").Value = 3749.8333
$ws.Range("K38").Value = 4084.6665
$ws.Range("L38").Value = 11249.4999
$ws.Range("M38").Value = -3712.6665
$ws.Range("N38").Value = -11993.4999
$ws.Range("H41").Value = 1393.8889
$ws.Range("I41").Value = 387.25
$ws.Range("J41").Value = 2199.2
$ws.Range("K41").Value = 387.25
$ws.Range("L41").Value = 2199.2
$ws.Range("M41").Value = 52.75
$ws.Range("N41").Value = -3079.2
$ws.Range("H53").Value = 5537.375
$ws.Range("I53").Value = 6042.857
$ws.Range("K53").Value = 6042.857
$ws.Range("M53").Value = -5405.857
$ws.Range("H64").Value = 50538.668
$ws.Range("I64").Value = 62549.715
$ws.Range("K64").Value = 62549.715
$ws.Range("M64").Value = -62301.715
$ws.Range("H67").Value = 50538.668
$ws.Range("I67").Value = 62549.715
$ws.Range("K67").Value = 62549.715
$ws.Range("M67").Value = -61691.715
$ws.Range("H98").Value = 32896.05
$ws.Range("I98").Value = 33113.293
$ws.Range("J98").Value = 31665
$ws.Range("K98").Value = 33113.293
$ws.Range("L98").Value = 31665
$ws.Range("M98").Value = -31615.293
$ws.Range("N98").Value = -34661
$ws.Range("H116").Value = 622493.9399999999
$ws.Range("I116").Value = 1239166.5
$ws.Range("J116").Value = 5821.3335
$ws.Range("K116").Value = 1239166.5
$ws.Range("L116").Value = 5821.3335
$ws.Range("M116").Value = -1235724.5
$ws.Range("N116").Value = -12705.3335
$ws.Range("H122").Value = 32896.05
$ws.Range("I122").Value = 33113.293
$ws.Range("J122").Value = 31665
$ws.Range("K122").Value = 99339.87899999999
$ws.Range("L122").Value = 94995
$ws.Range("M122").Value = -96889.87899999999
$ws.Range("N122").Value = -99895
$ws.Range("H138").Value = 3945.3845
$ws.Range("J138").Value = 4454.4688
$ws.Range("L138").Value = 13363.4064
$ws.Range("N138").Value = -23643.4064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2834.6978
$ws.Range("I32").Value = 2925.4
$ws.Range("J32").Value = 1625.3334
$ws.Range("K32").Value = 2925.4
$ws.Range("L32").Value = 1625.3334
$ws.Range("M32").Value = -2638.4
$ws.Range("N32").Value = -2199.3334
$ws.Range("H45").Value = 4394.815
$ws.Range("I45").Value = 3855.647
$ws.Range("J45").Value = 5311.4
$ws.Range("K45").Value = 3855.647
$ws.Range("L45").Value = 5311.4
$ws.Range("M45").Value = -3478.647
$ws.Range("N45").Value = -6065.4
$ws.Range("H74").Value = 2529.158
$ws.Range("I74").Value = 1375.7872
$ws.Range("K74").Value = 1375.7872
$ws.Range("M74").Value = -501.7872
$ws.Range("H77").Value = 2529.158
$ws.Range("I77").Value = 1375.7872
$ws.Range("K77").Value = 6878.936
$ws.Range("M77").Value = -2510.936
$ws.Range("H95").Value = 62999.668
$ws.Range("J95").Value = 62999.668
$ws.Range("L95").Value = 62999.668
$ws.Range("N95").Value = -68491.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 75167.5
$ws.Range("J122").Value = 75167.5
$ws.Range("L122").Value = 75167.5
$ws.Range("N122").Value = -84967.5
$ws.Range("H134").Value = 2763.838
$ws.Range("I134").Value = 2403.0356
$ws.Range("J134").Value = 3886.3333
$ws.Range("K134").Value = 7209.1068
$ws.Range("L134").Value = 11658.9999
$ws.Range("M134").Value = -4674.1068
$ws.Range("N134").Value = -16728.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H16").Value = 1201.6
$ws.Range("I16").Value = 1201.6
$ws.Range("K16").Value = 1201.6
$ws.Range("M16").Value = -914.5999999999999
$ws.Range("H19").Value = 481.30768
$ws.Range("I19").Value = 485.7
$ws.Range("J19").Value = 466.66666
$ws.Range("K19").Value = 485.7
$ws.Range("L19").Value = 466.66666
$ws.Range("M19").Value = -315.7
$ws.Range("N19").Value = -806.66666
$ws.Range("H24").Value = 481.30768
$ws.Range("I24").Value = 485.7
$ws.Range("J24").Value = 466.66666
$ws.Range("K24").Value = 485.7
$ws.Range("L24").Value = 466.66666
$ws.Range("M24").Value = -315.7
$ws.Range("N24").Value = -806.66666
$ws.Range("H31").Value = 2705.5806
$ws.Range("J31").Value = 4762.8184
$ws.Range("L31").Value = 4762.8184
$ws.Range("N31").Value = -5352.8184
$ws.Range("H34").Value = 2705.5806
$ws.Range("J34").Value = 4762.8184
$ws.Range("L34").Value = 4762.8184
$ws.Range("N34").Value = -5166.8184
$ws.Range("H59").Value = 49999.6
$ws.Range("J59").Value = 49999.6
$ws.Range("L59").Value = 49999.6
$ws.Range("N59").Value = -52289.6
$ws.Range("H113").Value = 1201.6
$ws.Range("I113").Value = 1201.6
$ws.Range("K113").Value = 1201.6
$ws.Range("M113").Value = 968.4000000000001
$ws.Range("H132").Value = 20673.625
$ws.Range("I132").Value = 4250.3687
$ws.Range("K132").Value = 12751.1061
$ws.Range("M132").Value = -10221.1061
$ws.Range("H141").Value = 82809.5
$ws.Range("I141").Value = 49749.5
$ws.Range("K141").Value = 49749.5
$ws.Range("M141").Value = -44569.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 6812.7144
$ws.Range("I57").Value = 3137.8
$ws.Range("K57").Value = 9413.400000000001
$ws.Range("M57").Value = -8854.400000000001
$ws.Range("H97").Value = 27861.088
$ws.Range("I97").Value = 39151.812
$ws.Range("J97").Value = 2053.7144
$ws.Range("K97").Value = 117455.436
$ws.Range("L97").Value = 6161.1432
$ws.Range("M97").Value = -116959.436
$ws.Range("N97").Value = -7153.1432
$ws.Range("H115").Value = 9593
$ws.Range("I115").Value = 9593
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 28779
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -27604
$ws.Range("N115").ClearContents()
$ws.Range("H126").Value = 13789.889
$ws.Range("I126").Value = 2259.6667
$ws.Range("K126").Value = 6779.000100000001
$ws.Range("M126").Value = -1839.000100000001
$ws.Range("H131").Value = 1875.2646
$ws.Range("J131").Value = 2175.7693
$ws.Range("L131").Value = 6527.3079
$ws.Range("N131").Value = -16607.3079
$ws.Range("H137").Value = 3642.8333
$ws.Range("J137").Value = 6413.3335
$ws.Range("L137").Value = 19240.0005
$ws.Range("N137").Value = -29440.0005
$ws.Range("H140").Value = 9311.360000000001
$ws.Range("I140").Value = 11376.889
$ws.Range("K140").Value = 34130.667
$ws.Range("M140").Value = -28950.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 5109.486
$ws.Range("I102").Value = 5822.8213
$ws.Range("K102").Value = 5822.8213
$ws.Range("M102").Value = -4200.8213
$ws.Range("H126").Value = 16722.174
$ws.Range("I126").Value = 19543.3
$ws.Range("J126").Value = 14552.077
$ws.Range("K126").Value = 58629.89999999999
$ws.Range("L126").Value = 43656.231
$ws.Range("M126").Value = -56159.89999999999
$ws.Range("N126").Value = -48596.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13213.091
$ws.Range("I22").Value = 18477.715
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 18477.715
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -18182.715
$ws.Range("N22").Value = -4590
$ws.Range("H27").Value = 13213.091
$ws.Range("I27").Value = 18477.715
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 18477.715
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -18370.715
$ws.Range("N27").Value = -4214
$ws.Range("H46").Value = 3544.111
$ws.Range("I46").Value = 871.2857
$ws.Range("K46").Value = 871.2857
$ws.Range("M46").Value = -683.2857
$ws.Range("H82").Value = 2730.5
$ws.Range("I82").Value = 4700.4
$ws.Range("K82").Value = 4700.4
$ws.Range("M82").Value = -4339.4
$ws.Range("H85").Value = 2730.5
$ws.Range("I85").Value = 4700.4
$ws.Range("K85").Value = 4700.4
$ws.Range("M85").Value = -3452.4
$ws.Range("H132").Value = 624035.4399999999
$ws.Range("I132").Value = 1147688.6
$ws.Range("K132").Value = 3443065.8
$ws.Range("M132").Value = -3440535.8
